$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-06 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-07 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("42+31=73", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=81", 2) | Out-Null
$d.Content.Find.Execute("15-5=10", $true, $false, $false, $false, $false, $true, 1, $false, "44-12=32", 2) | Out-Null
$d.Content.Find.Execute("89-64=25", $true, $false, $false, $false, $false, $true, 1, $false, "31+11=42", 2) | Out-Null
$d.Content.Find.Execute("58-24=34", $true, $false, $false, $false, $false, $true, 1, $false, "64-49=15", 2) | Out-Null
$d.Content.Find.Execute("24+52=76", $true, $false, $false, $false, $false, $true, 1, $false, "20+44=64", 2) | Out-Null
$d.Content.Find.Execute("61-0=61", $true, $false, $false, $false, $false, $true, 1, $false, "41+42=83", 2) | Out-Null
$d.Content.Find.Execute("40-38=2", $true, $false, $false, $false, $false, $true, 1, $false, "94-84=10", 2) | Out-Null
$d.Content.Find.Execute("10+86=96", $true, $false, $false, $false, $false, $true, 1, $false, "71+8=79", 2) | Out-Null
$d.Content.Find.Execute("39+37=76", $true, $false, $false, $false, $false, $true, 1, $false, "30-12=18", 2) | Out-Null
$d.Content.Find.Execute("73+0=73", $true, $false, $false, $false, $false, $true, 1, $false, "51-42=9", 2) | Out-Null
$d.Content.Find.Execute("95-16=79", $true, $false, $false, $false, $false, $true, 1, $false, "34-2=32", 2) | Out-Null
$d.Content.Find.Execute("5+41=46", $true, $false, $false, $false, $false, $true, 1, $false, "29+12=41", 2) | Out-Null
$d.Content.Find.Execute("78-18=60", $true, $false, $false, $false, $false, $true, 1, $false, "29+13=42", 2) | Out-Null
$d.Content.Find.Execute("35+63=98", $true, $false, $false, $false, $false, $true, 1, $false, "98-38=60", 2) | Out-Null
$d.Content.Find.Execute("89-25=64", $true, $false, $false, $false, $false, $true, 1, $false, "1+20=21", 2) | Out-Null
$d.Content.Find.Execute("98-7=91", $true, $false, $false, $false, $false, $true, 1, $false, "73+1=74", 2) | Out-Null
$d.Content.Find.Execute("24-8=16", $true, $false, $false, $false, $false, $true, 1, $false, "70-65=5", 2) | Out-Null
$d.Content.Find.Execute("64-2=62", $true, $false, $false, $false, $false, $true, 1, $false, "86-51=35", 2) | Out-Null
$d.Content.Find.Execute("13+73=86", $true, $false, $false, $false, $false, $true, 1, $false, "61-45=16", 2) | Out-Null
$d.Content.Find.Execute("80-28=52", $true, $false, $false, $false, $false, $true, 1, $false, "34+39=73", 2) | Out-Null
$d.Content.Find.Execute("84-16=68", $true, $false, $false, $false, $false, $true, 1, $false, "43+56=99", 2) | Out-Null
$d.Content.Find.Execute("0+53=53", $true, $false, $false, $false, $false, $true, 1, $false, "74-56=18", 2) | Out-Null
$d.Content.Find.Execute("29-9=20", $true, $false, $false, $false, $false, $true, 1, $false, "32-7=25", 2) | Out-Null
$d.Content.Find.Execute("29+60=89", $true, $false, $false, $false, $false, $true, 1, $false, "9+72=81", 2) | Out-Null
$d.Content.Find.Execute("10+88=98", $true, $false, $false, $false, $false, $true, 1, $false, "26+68=94", 2) | Out-Null
$d.Content.Find.Execute("17+36=53", $true, $false, $false, $false, $false, $true, 1, $false, "57-15=42", 2) | Out-Null
$d.Content.Find.Execute("60+1=61", $true, $false, $false, $false, $false, $true, 1, $false, "69-3=66", 2) | Out-Null
$d.Content.Find.Execute("94-11=83", $true, $false, $false, $false, $false, $true, 1, $false, "82-0=82", 2) | Out-Null
$d.Content.Find.Execute("80+15=95", $true, $false, $false, $false, $false, $true, 1, $false, "4+54=58", 2) | Out-Null
$d.Content.Find.Execute("40-14=26", $true, $false, $false, $false, $false, $true, 1, $false, "93+5=98", 2) | Out-Null
$d.Content.Find.Execute("44+47=91", $true, $false, $false, $false, $false, $true, 1, $false, "6-0=6", 2) | Out-Null
$d.Content.Find.Execute("62-28=34", $true, $false, $false, $false, $false, $true, 1, $false, "56-18=38", 2) | Out-Null
$d.Content.Find.Execute("71+17=88", $true, $false, $false, $false, $false, $true, 1, $false, "65-41=24", 2) | Out-Null
$d.Content.Find.Execute("11+52=63", $true, $false, $false, $false, $false, $true, 1, $false, "39-4=35", 2) | Out-Null
$d.Content.Find.Execute("79-33=46", $true, $false, $false, $false, $false, $true, 1, $false, "85-9=76", 2) | Out-Null
$d.Content.Find.Execute("24+56=80", $true, $false, $false, $false, $false, $true, 1, $false, "20+73=93", 2) | Out-Null
$d.Content.Find.Execute("34-18=16", $true, $false, $false, $false, $false, $true, 1, $false, "75+15=90", 2) | Out-Null
$d.Content.Find.Execute("66-48=18", $true, $false, $false, $false, $false, $true, 1, $false, "45-34=11", 2) | Out-Null
$d.Content.Find.Execute("6+21=27", $true, $false, $false, $false, $false, $true, 1, $false, "3+86=89", 2) | Out-Null
$d.Content.Find.Execute("95-32=63", $true, $false, $false, $false, $false, $true, 1, $false, "6-3=3", 2) | Out-Null
$d.Content.Find.Execute("28-13=15", $true, $false, $false, $false, $false, $true, 1, $false, "47+7=54", 2) | Out-Null
$d.Content.Find.Execute("18-6=12", $true, $false, $false, $false, $false, $true, 1, $false, "67-22=45", 2) | Out-Null
$d.Content.Find.Execute("30+67=97", $true, $false, $false, $false, $false, $true, 1, $false, "61-5=56", 2) | Out-Null
$d.Content.Find.Execute("65-46=19", $true, $false, $false, $false, $false, $true, 1, $false, "69+24=93", 2) | Out-Null
$d.Content.Find.Execute("87-87=0", $true, $false, $false, $false, $false, $true, 1, $false, "94-38=56", 2) | Out-Null
$d.Content.Find.Execute("73-44=29", $true, $false, $false, $false, $false, $true, 1, $false, "51+14=65", 2) | Out-Null
$d.Content.Find.Execute("8+81=89", $true, $false, $false, $false, $false, $true, 1, $false, "57+21=78", 2) | Out-Null
$d.Content.Find.Execute("99-59=40", $true, $false, $false, $false, $false, $true, 1, $false, "37+32=69", 2) | Out-Null
$d.Content.Find.Execute("3+41=44", $true, $false, $false, $false, $false, $true, 1, $false, "51-30=21", 2) | Out-Null
$d.Content.Find.Execute("72-42=30", $true, $false, $false, $false, $false, $true, 1, $false, "25+27=52", 2) | Out-Null
$d.Content.Find.Execute("95-13=82", $true, $false, $false, $false, $false, $true, 1, $false, "81+14=95", 2) | Out-Null
$d.Content.Find.Execute("61-7=54", $true, $false, $false, $false, $false, $true, 1, $false, "23+53=76", 2) | Out-Null
$d.Content.Find.Execute("74-33=41", $true, $false, $false, $false, $false, $true, 1, $false, "67-14=53", 2) | Out-Null
$d.Content.Find.Execute("91-34=57", $true, $false, $false, $false, $false, $true, 1, $false, "50+14=64", 2) | Out-Null
$d.Content.Find.Execute("46-27=19", $true, $false, $false, $false, $false, $true, 1, $false, "42-33=9", 2) | Out-Null
$d.Content.Find.Execute("76-47=29", $true, $false, $false, $false, $false, $true, 1, $false, "40-39=1", 2) | Out-Null
$d.Content.Find.Execute("21+39=60", $true, $false, $false, $false, $false, $true, 1, $false, "70-66=4", 2) | Out-Null
$d.Content.Find.Execute("86-48=38", $true, $false, $false, $false, $false, $true, 1, $false, "8+58=66", 2) | Out-Null
$d.Content.Find.Execute("51+20=71", $true, $false, $false, $false, $false, $true, 1, $false, "94-50=44", 2) | Out-Null
$d.Content.Find.Execute("96-72=24", $true, $false, $false, $false, $false, $true, 1, $false, "7+46=53", 2) | Out-Null
$d.Content.Find.Execute("62+23=85", $true, $false, $false, $false, $false, $true, 1, $false, "55-31=24", 2) | Out-Null
$d.Content.Find.Execute("74-24=50", $true, $false, $false, $false, $false, $true, 1, $false, "50+15=65", 2) | Out-Null
$d.Content.Find.Execute("62-38=24", $true, $false, $false, $false, $false, $true, 1, $false, "75-1=74", 2) | Out-Null
$d.Content.Find.Execute("27+37=64", $true, $false, $false, $false, $false, $true, 1, $false, "61-42=19", 2) | Out-Null
$d.Content.Find.Execute("12+74=86", $true, $false, $false, $false, $false, $true, 1, $false, "5+42=47", 2) | Out-Null
$d.Content.Find.Execute("37-20=17", $true, $false, $false, $false, $false, $true, 1, $false, "33+4=37", 2) | Out-Null
$d.Content.Find.Execute("89-82=7", $true, $false, $false, $false, $false, $true, 1, $false, "82-9=73", 2) | Out-Null
$d.Content.Find.Execute("43-28=15", $true, $false, $false, $false, $false, $true, 1, $false, "76-9=67", 2) | Out-Null
$d.Content.Find.Execute("72+11=83", $true, $false, $false, $false, $false, $true, 1, $false, "40+22=62", 2) | Out-Null
$d.Content.Find.Execute("70-27=43", $true, $false, $false, $false, $false, $true, 1, $false, "30+41=71", 2) | Out-Null
$d.Content.Find.Execute("22-13=9", $true, $false, $false, $false, $false, $true, 1, $false, "99-25=74", 2) | Out-Null
$d.Content.Find.Execute("89-19=70", $true, $false, $false, $false, $false, $true, 1, $false, "6+88=94", 2) | Out-Null
$d.Content.Find.Execute("47-32=15", $true, $false, $false, $false, $false, $true, 1, $false, "8+48=56", 2) | Out-Null
$d.Content.Find.Execute("20+35=55", $true, $false, $false, $false, $false, $true, 1, $false, "35+11=46", 2) | Out-Null
$d.Content.Find.Execute("47-17=30", $true, $false, $false, $false, $false, $true, 1, $false, "14+77=91", 2) | Out-Null
$d.Content.Find.Execute("30-21=9", $true, $false, $false, $false, $false, $true, 1, $false, "56+18=74", 2) | Out-Null
$d.Content.Find.Execute("37-14=23", $true, $false, $false, $false, $false, $true, 1, $false, "85-1=84", 2) | Out-Null
$d.Content.Find.Execute("92-76=16", $true, $false, $false, $false, $false, $true, 1, $false, "99-75=24", 2) | Out-Null
$d.Content.Find.Execute("58+3=61", $true, $false, $false, $false, $false, $true, 1, $false, "52-52=0", 2) | Out-Null
$d.Content.Find.Execute("96-42=54", $true, $false, $false, $false, $false, $true, 1, $false, "17-15=2", 2) | Out-Null
$d.Content.Find.Execute("51+40=91", $true, $false, $false, $false, $false, $true, 1, $false, "87-76=11", 2) | Out-Null
$d.Content.Find.Execute("8+35=43", $true, $false, $false, $false, $false, $true, 1, $false, "99-95=4", 2) | Out-Null
$d.Content.Find.Execute("43+42=85", $true, $false, $false, $false, $false, $true, 1, $false, "75-22=53", 2) | Out-Null
$d.Content.Find.Execute("83-51=32", $true, $false, $false, $false, $false, $true, 1, $false, "26+51=77", 2) | Out-Null
$d.Content.Find.Execute("6+58=64", $true, $false, $false, $false, $false, $true, 1, $false, "27-9=18", 2) | Out-Null
$d.Content.Find.Execute("76-41=35", $true, $false, $false, $false, $false, $true, 1, $false, "78-71=7", 2) | Out-Null
$d.Content.Find.Execute("92+2=94", $true, $false, $false, $false, $false, $true, 1, $false, "32+10=42", 2) | Out-Null
$d.Content.Find.Execute("75-12=63", $true, $false, $false, $false, $false, $true, 1, $false, "11+1=12", 2) | Out-Null
$d.Content.Find.Execute("78+9=87", $true, $false, $false, $false, $false, $true, 1, $false, "94-29=65", 2) | Out-Null
$d.Content.Find.Execute("49-14=35", $true, $false, $false, $false, $false, $true, 1, $false, "44-30=14", 2) | Out-Null
$d.Content.Find.Execute("70+2=72", $true, $false, $false, $false, $false, $true, 1, $false, "97-8=89", 2) | Out-Null
$d.Content.Find.Execute("11+48=59", $true, $false, $false, $false, $false, $true, 1, $false, "45-45=0", 2) | Out-Null
$d.Content.Find.Execute("17+71=88", $true, $false, $false, $false, $false, $true, 1, $false, "5+87=92", 2) | Out-Null
$d.Content.Find.Execute("67-46=21", $true, $false, $false, $false, $false, $true, 1, $false, "45+54=99", 2) | Out-Null
$d.Content.Find.Execute("61-22=39", $true, $false, $false, $false, $false, $true, 1, $false, "92-75=17", 2) | Out-Null
$d.Content.Find.Execute("36+55=91", $true, $false, $false, $false, $false, $true, 1, $false, "48-41=7", 2) | Out-Null
$d.Content.Find.Execute("56-53=3", $true, $false, $false, $false, $false, $true, 1, $false, "33+35=68", 2) | Out-Null
$d.Content.Find.Execute("7+8=15", $true, $false, $false, $false, $false, $true, 1, $false, "64-13=51", 2) | Out-Null
$d.Content.Find.Execute("68-36=32", $true, $false, $false, $false, $false, $true, 1, $false, "99-39=60", 2) | Out-Null
$d.Content.Find.Execute("7-5=2", $true, $false, $false, $false, $false, $true, 1, $false, "45-30=15", 2) | Out-Null
